$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing rows 189:197 down to 192:200, making room for 3 new
# rows of data at 189:191 (weekly update pattern: new readings inserted at
# top of this market/product block, older ones pushed down).
$ws.Range("A189:T191").Insert()

# ---- Row 189 ------------------------------------------------------------
$ws.Cells.Item(189, 1).Value = 3
$ws.Cells.Item(189, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(189, 3).Value = "Coquimbo"
$ws.Cells.Item(189, 4).Value = 45008
$ws.Cells.Item(189, 5).Value = 5
$ws.Cells.Item(189, 6).Value = "Fruta"
$ws.Cells.Item(189, 7).Value = 100107
$ws.Cells.Item(189, 8).Value = "Otros"
$ws.Cells.Item(189, 9).Value = 100107011
$ws.Cells.Item(189, 10).Value = "Tuna"
$ws.Cells.Item(189, 11).Value = "Sin especificar"
$ws.Cells.Item(189, 12).Value = "Especial"
$ws.Cells.Item(189, 13).Value = 68
$ws.Cells.Item(189, 14).Value = 16000
$ws.Cells.Item(189, 15).Value = 16000
$ws.Cells.Item(189, 16).Value = 16000
$ws.Cells.Item(189, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(189, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(189, 19).Value = 1000
$ws.Cells.Item(189, 20).Value = 16

# ---- Row 190 ------------------------------------------------------------
$ws.Cells.Item(190, 1).Value = 3
$ws.Cells.Item(190, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(190, 3).Value = "Coquimbo"
$ws.Cells.Item(190, 4).Value = 45008
$ws.Cells.Item(190, 5).Value = 5
$ws.Cells.Item(190, 6).Value = "Fruta"
$ws.Cells.Item(190, 7).Value = 100107
$ws.Cells.Item(190, 8).Value = "Otros"
$ws.Cells.Item(190, 9).Value = 100107011
$ws.Cells.Item(190, 10).Value = "Tuna"
$ws.Cells.Item(190, 11).Value = "Sin especificar"
$ws.Cells.Item(190, 12).Value = "Primera"
$ws.Cells.Item(190, 13).Value = 67
$ws.Cells.Item(190, 14).Value = 13000
$ws.Cells.Item(190, 15).Value = 13000
$ws.Cells.Item(190, 16).Value = 13000
$ws.Cells.Item(190, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(190, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(190, 19).Value = 812
$ws.Cells.Item(190, 20).Value = 16

# ---- Row 191 ------------------------------------------------------------
$ws.Cells.Item(191, 1).Value = 3
$ws.Cells.Item(191, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(191, 3).Value = "Coquimbo"
$ws.Cells.Item(191, 4).Value = 45008
$ws.Cells.Item(191, 5).Value = 5
$ws.Cells.Item(191, 6).Value = "Fruta"
$ws.Cells.Item(191, 7).Value = 100107
$ws.Cells.Item(191, 8).Value = "Otros"
$ws.Cells.Item(191, 9).Value = 100107011
$ws.Cells.Item(191, 10).Value = "Tuna"
$ws.Cells.Item(191, 11).Value = "Sin especificar"
$ws.Cells.Item(191, 12).Value = "Segunda"
$ws.Cells.Item(191, 13).Value = 60
$ws.Cells.Item(191, 14).Value = 10000
$ws.Cells.Item(191, 15).Value = 10000
$ws.Cells.Item(191, 16).Value = 10000
$ws.Cells.Item(191, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(191, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(191, 19).Value = 625
$ws.Cells.Item(191, 20).Value = 16

$wb.Save()
